$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / sheet name
$ws.Name = "linear+or"

# New data values (row -> column A..K), reflecting the reshuffled / updated
# results (hiv ba results).
$data = @{
    2 = @("abuse", "69.668 [36.641 -- 89.251]", "30.332 [10.749 -- 63.359]", "0.810 [0.243 -- 2.105]", "0.353 [0.141 -- 0.750]", "1.163 [0.584 -- 2.630]", "1.315 [0.657 -- 2.159]", "1.297 [0.625 -- 2.103]", "0.647 [0.596 -- 0.754]", "0.639 [0.628 -- 0.672]", "-10.067 [-10.346 -- -9.734]")
    3 = @("heart disorder", "77.584 [52.898 -- 94.606]", "22.416 [5.394 -- 47.102]", "0.903 [0.342 -- 2.384]", "0.261 [0.067 -- 0.526]", "1.164 [0.585 -- 2.632]", "1.129 [0.547 -- 1.965]", "1.126 [0.520 -- 1.950]", "0.647 [0.596 -- 0.754]", "0.786 [0.771 -- 0.803]", "-10.067 [-10.346 -- -9.734]")
    4 = @("avg", "80.314 [61.774 -- 89.459]", "19.686 [10.541 -- 38.226]", "0.935 [0.391 -- 2.293]", "0.229 [0.127 -- 0.424]", "1.164 [0.585 -- 2.631]", "0.959 [0.578 -- 1.430]", "0.952 [0.573 -- 1.396]", "0.647 [0.596 -- 0.754]", "nan [nan -- nan]", "-10.067 [-10.346 -- -9.734]")
    5 = @("sleep disorder", "93.686 [71.436 -- 103.977]", "6.314 [-3.977 -- 28.564]", "1.090 [0.480 -- 2.426]", "0.073 [-0.048 -- 0.371]", "1.164 [0.585 -- 2.632]", "0.433 [-0.145 -- 1.090]", "0.433 [-0.145 -- 1.090]", "0.647 [0.596 -- 0.754]", "0.641 [0.626 -- 0.666]", "-10.067 [-10.346 -- -9.734]")
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
